$d = $word.ActiveDocument

foreach ($p in $d.Paragraphs) {
    if ($p.Range.Text -like "Docente(s)*") {
        $r = $p.Range
        $r.InsertParagraphAfter()
        $newPara = $p.Next()
        $npr = $newPara.Range

        $xml = '<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">' +
               '<w:pPr><w:pStyle w:val="ListBullet"/></w:pPr>' +
               '<w:r><w:t>7459752 - Maria Ismenia Sodero Toledo Faria</w:t><w:br/></w:r>' +
               '<w:r><w:t>2166002 - Sandra Giacomin Schneider</w:t><w:br/></w:r>' +
               '<w:r><w:t>1922320 - Sebastiao Ribeiro</w:t></w:r>' +
               '</w:p>'

        $npr.InsertXML($xml)
        break
    }
}
